$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Content fixes -------------------------------------------------------
# Fix typo in the Features text of row 2 ("inagenes" -> "imagenes")
$ws.Range("B2").Value = "Captura de imagenes y procesamiento inicial de imagenes"

# Replace the Features text of row 15 with the new wording
$ws.Range("B15").Value = "relacion entre la interaccion de video con acciones de vision"

# --- Cosmetic cleanup ------------------------------------------------------
# Collapse the stray fully-empty spacer rows (their explicit row height goes
# away and they disappear from the saved XML once no attribute differs from
# the sheet default) and drop the custom row height Excel had stamped on
# every data row.
$ws.Range("A1:F27").EntireRow.AutoFit()

# Row 2 is taller on purpose (it wraps to two lines) - restore that custom
# height since AutoFit reset it above.
$ws.Rows("2:2").RowHeight = 22.35

# Reflect the author's final cursor/selection position in the saved view.
$null = $ws.Range("B26:B27").Select()
